$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the "Chirimoya" block at
# Vega Modelo de Temuco. Insert a fresh row above the current row 251 so
# that all the existing rows (251-267) shift down by one (to 252-268),
# then populate the newly inserted row with the new week's data.
$ws.Rows.Item(251).Insert()

$row = 251
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 45267
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100107
$ws.Cells.Item($row, 8).Value = "Otros"
$ws.Cells.Item($row, 9).Value = 100107002
$ws.Cells.Item($row, 10).Value = "Chirimoya"
$ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 85
$ws.Cells.Item($row, 14).Value = 2300
$ws.Cells.Item($row, 15).Value = 2300
$ws.Cells.Item($row, 16).Value = 2300
$ws.Cells.Item($row, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item($row, 18).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 19).Value = 2300
$ws.Cells.Item($row, 20).Value = 1
